$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 needs to hold the *text* "1" (not the number 1) while keeping its
# existing cell formatting (style stays the same "General" xf it already
# has). A direct `Range.Value = "1"` would be auto-coerced to the number
# 1 by Excel's usual type inference, and NumberFormat="@" (or a leading
# apostrophe) would stamp a different/quote-prefixed style onto the cell.
# So: stage the text value in a scratch cell formatted as Text, copy it,
# and paste-special *values only* into B11 - that carries over the text
# type without disturbing B11's own formatting.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
